$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13711
$ws.Range("E2").Value = 950
$ws.Range("F2").Value = 950
$ws.Range("G2").Value = 948
$ws.Range("H2").Value = 733
$ws.Range("I2").Value = 726
$ws.Range("J2").Value = 7
$ws.Range("K2").Value = 7216
$ws.Range("L2").Value = 2068
$ws.Range("M2").Value = 5148
$ws.Range("N2").Value = 5050
$ws.Range("O2").Value = 98
$ws.Range("P2").Value = 295
$ws.Range("Q2").Value = 502
$ws.Range("R2").Value = -902
$ws.Range("S2").Value = 123
$ws.Range("T2").Value = 459
$ws.Range("U2").Value = 43
$ws.Range("V2").Value = 624
$ws.Range("W2").Value = 6.93
$ws.Range("X2").Value = 5.34
$ws.Range("Y2").Value = 15.31
$ws.Range("Z2").Value = 10.84
$ws.Range("AA2").Value = 40.17
$ws.Range("AB2").Value = 1616.87
$ws.Range("AC2").Value = 1230
$ws.Range("AD2").Value = 11.39
$ws.Range("AE2").Value = 8549
$ws.Range("AF2").Value = 1.64
$ws.Range("AG2").Value = 175
$ws.Range("AH2").Value = 1.25
$ws.Range("AI2").Value = 14.23
$ws.Range("AJ2").Value = 59070000
$ws.Range("D3").Value = 12212
$ws.Range("E3").Value = 560
$ws.Range("F3").Value = 560
$ws.Range("G3").Value = 462
$ws.Range("H3").Value = 322
$ws.Range("I3").Value = 346
$ws.Range("J3").Value = -24
$ws.Range("K3").Value = 7175
$ws.Range("L3").Value = 1793
$ws.Range("M3").Value = 5382
$ws.Range("N3").Value = 5302
$ws.Range("O3").Value = 80
$ws.Range("P3").Value = 295
$ws.Range("Q3").Value = 616
$ws.Range("R3").Value = 586
$ws.Range("S3").Value = -173
$ws.Range("T3").Value = 286
$ws.Range("U3").Value = 330
$ws.Range("V3").Value = 595
$ws.Range("W3").Value = 4.59
$ws.Range("X3").Value = 2.63
$ws.Range("Y3").Value = 6.68
$ws.Range("Z3").Value = 4.47
$ws.Range("AA3").Value = 33.31
$ws.Range("AB3").Value = 1701.84
$ws.Range("AC3").Value = 585
$ws.Range("AD3").Value = 23.4
$ws.Range("AE3").Value = 8976
$ws.Range("AF3").Value = 1.53
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 1.46
$ws.Range("AI3").Value = 34.16
$ws.Range("AJ3").Value = 59070000
$ws.Range("D4").Value = 11177
$ws.Range("E4").Value = 853
$ws.Range("F4").Value = 853
$ws.Range("G4").Value = 648
$ws.Range("H4").Value = 445
$ws.Range("I4").Value = 437
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 7239
$ws.Range("L4").Value = 1545
$ws.Range("M4").Value = 5694
$ws.Range("N4").Value = 5625
$ws.Range("O4").Value = 69
$ws.Range("P4").Value = 295
$ws.Range("Q4").Value = 851
$ws.Range("R4").Value = -1063
$ws.Range("S4").Value = -223
$ws.Range("T4").Value = 233
$ws.Range("U4").Value = 618
$ws.Range("V4").Value = 370
$ws.Range("W4").Value = 7.63
$ws.Range("X4").Value = 3.98
$ws.Range("Y4").Value = 8
$ws.Range("Z4").Value = 6.18
$ws.Range("AA4").Value = 27.13
$ws.Range("AB4").Value = 1806.91
$ws.Range("AC4").Value = 740
$ws.Range("AD4").Value = 16.22
$ws.Range("AE4").Value = 9523
$ws.Range("AF4").Value = 1.26
$ws.Range("AG4").Value = 300
$ws.Range("AH4").Value = 2.5
$ws.Range("AI4").Value = 40.56
$ws.Range("AJ4").Value = 59070000
$ws.Range("D5").Value = 11972
$ws.Range("E5").Value = 1040
$ws.Range("F5").Value = 1040
$ws.Range("G5").Value = 1256
$ws.Range("H5").Value = 1040
$ws.Range("I5").Value = 1036
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 8358
$ws.Range("L5").Value = 1836
$ws.Range("M5").Value = 6522
$ws.Range("N5").Value = 6445
$ws.Range("O5").Value = 77
$ws.Range("P5").Value = 295
$ws.Range("Q5").Value = 555
$ws.Range("R5").Value = -462
$ws.Range("S5").Value = -230
$ws.Range("T5").Value = 594
$ws.Range("U5").Value = -39
$ws.Range("V5").Value = 278
$ws.Range("W5").Value = 8.68
$ws.Range("X5").Value = 8.68
$ws.Range("Y5").Value = 17.16
$ws.Range("Z5").Value = 13.33
$ws.Range("AA5").Value = 28.16
$ws.Range("AB5").Value = 2094.19
$ws.Range("AC5").Value = 1753
$ws.Range("AD5").Value = 22.7
$ws.Range("AE5").Value = 10910
$ws.Range("AF5").Value = 3.65
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 0.88
$ws.Range("AI5").Value = 19.96
$ws.Range("AJ5").Value = 59070000
$ws.Range("D6").Value = 13836
$ws.Range("E6").Value = 1063
$ws.Range("F6").Value = 1063
$ws.Range("G6").Value = 1766
$ws.Range("H6").Value = 1328
$ws.Range("I6").Value = 1322
$ws.Range("K6").Value = 9476
$ws.Range("L6").Value = 1895
$ws.Range("M6").Value = 7581
$ws.Range("N6").Value = 7496
$ws.Range("P6").Value = 295
$ws.Range("Q6").Value = 846
$ws.Range("R6").Value = -387
$ws.Range("S6").Value = -294
$ws.Range("T6").Value = 1009
$ws.Range("U6").Value = -163
$ws.Range("V6").Value = 202
$ws.Range("W6").Value = 7.68
$ws.Range("X6").Value = 9.59
$ws.Range("Y6").Value = 18.96
$ws.Range("Z6").Value = 14.89
$ws.Range("AA6").Value = 25
$ws.Range("AB6").Value = 2446.59
$ws.Range("AC6").Value = 2237
$ws.Range("AD6").Value = 28.47
$ws.Range("AE6").Value = 12690
$ws.Range("AF6").Value = 5.02
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 0.63
$ws.Range("AI6").Value = 17.88
$ws.Range("AJ6").Value = 59070000
$ws.Range("D7").Value = 14852
$ws.Range("E7").Value = 918
$ws.Range("G7").Value = 1274
$ws.Range("H7").Value = 1058
$ws.Range("I7").Value = 1089
$ws.Range("K7").Value = 13776
$ws.Range("L7").Value = 4657
$ws.Range("M7").Value = 9117
$ws.Range("N7").Value = 9026
$ws.Range("P7").Value = 301
$ws.Range("Q7").Value = 956
$ws.Range("R7").Value = -3120
$ws.Range("S7").Value = 1954
$ws.Range("T7").Value = 3592
$ws.Range("U7").Value = -2102
$ws.Range("W7").Value = 6.18
$ws.Range("X7").Value = 7.12
$ws.Range("Y7").Value = 13.18
$ws.Range("Z7").Value = 9.1
$ws.Range("AA7").Value = 51.08
$ws.Range("AC7").Value = 1802
$ws.Range("AD7").Value = 33.68
$ws.Range("AE7").Value = 14799
$ws.Range("AF7").Value = 4.1
$ws.Range("AG7").Value = 390
$ws.Range("AH7").Value = 0.64
$ws.Range("AI7").Value = 21.84
$ws.Range("D8").Value = 19676
$ws.Range("E8").Value = 1318
$ws.Range("G8").Value = 1519
$ws.Range("H8").Value = 1206
$ws.Range("I8").Value = 1354
$ws.Range("K8").Value = 17283
$ws.Range("L8").Value = 6849
$ws.Range("M8").Value = 10433
$ws.Range("N8").Value = 10308
$ws.Range("P8").Value = 301
$ws.Range("Q8").Value = 1378
$ws.Range("R8").Value = -2918
$ws.Range("S8").Value = 1735
$ws.Range("T8").Value = 2750
$ws.Range("U8").Value = -574
$ws.Range("W8").Value = 6.7
$ws.Range("X8").Value = 6.13
$ws.Range("Y8").Value = 13.89
$ws.Range("Z8").Value = 7.77
$ws.Range("AA8").Value = 65.65000000000001
$ws.Range("AC8").Value = 2220
$ws.Range("AD8").Value = 25.86
$ws.Range("AE8").Value = 16901
$ws.Range("AF8").Value = 3.4
$ws.Range("AG8").Value = 394
$ws.Range("AH8").Value = 0.6899999999999999
$ws.Range("AI8").Value = 17.74
$ws.Range("D9").Value = 25051
$ws.Range("E9").Value = 1750
$ws.Range("G9").Value = 1967
$ws.Range("H9").Value = 1558
$ws.Range("I9").Value = 1772
$ws.Range("K9").Value = 19891
$ws.Range("L9").Value = 7909
$ws.Range("M9").Value = 11982
$ws.Range("N9").Value = 11852
$ws.Range("P9").Value = 301
$ws.Range("Q9").Value = 1807
$ws.Range("R9").Value = -2243
$ws.Range("S9").Value = 544
$ws.Range("T9").Value = 2041
$ws.Range("U9").Value = -254
$ws.Range("W9").Value = 6.99
$ws.Range("X9").Value = 6.22
$ws.Range("Y9").Value = 15.99
$ws.Range("Z9").Value = 8.380000000000001
$ws.Range("AA9").Value = 66.01000000000001
$ws.Range("AC9").Value = 2906
$ws.Range("AD9").Value = 19.75
$ws.Range("AE9").Value = 19433
$ws.Range("AF9").Value = 2.95
$ws.Range("AG9").Value = 394
$ws.Range("AH9").Value = 0.6899999999999999
$ws.Range("AI9").Value = 13.55
